# Adds a new requirement/use-case row (row 19) to the "Requisitos e Casos de
# Uso" sheet: REQ# 15 / "Selecionar Bebê" use case #14, describing that the
# app should ask which baby to view when more than one is registered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - sequential requirement number
$ws.Cells.Item(19, 1).Value = 15

# Column B - requirement description (wrap-text style, like similarly
# multi-line rows e.g. B13/B14/B18)
$ws.Cells.Item(19, 2).Value = "Ao entrar no aplicativo, se existir mais de um bebê, perguntar ao usuário qual bebê deseja visualizar"

# Column C - sequential use-case number
$ws.Cells.Item(19, 3).Value = 14

# Column D - use case name
$ws.Cells.Item(19, 4).Value = "Selecionar Bebê"

# Match the row style used for other wrapped, taller description rows
# (row height 30, style index 3 = wrap text) applied to column B.
$ws.Range("B19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 30

# Keep the view/selection state in sync with the diff (scrolled down one
# row, active cell now on the new row's D column).
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("D20").Select()
